# Adds the Outbound* booking-flow columns (AL:BA) captured once a user
# reaches flight selection, and the three new booking rows (19-21) that
# were recorded with that additional detail - guardar informacion del vuelo.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns AL1:BA1 (row 1) ---
$ws.Range("AL1").Value = "Outbound Flight Date"
$ws.Range("AM1").Value = "Outbound Origin Code"
$ws.Range("AN1").Value = "Outbound Origin Name"
$ws.Range("AO1").Value = "Outbound Destination Code"
$ws.Range("AP1").Value = "Outbound Destination Name"
$ws.Range("AQ1").Value = "Outbound Arrival Next Day"
$ws.Range("AR1").Value = "Outbound Airline"
$ws.Range("AS1").Value = "Outbound Price"
$ws.Range("AT1").Value = "Outbound Direction"
$ws.Range("AU1").Value = "Outbound Selected Class"
$ws.Range("AV1").Value = "Outbound Selected Type"
$ws.Range("AW1").Value = "Outbound Cabin"
$ws.Range("AX1").Value = "Outbound Economy Option"
$ws.Range("AY1").Value = "Outbound Economy Price"
$ws.Range("AZ1").Value = "Outbound Premium Option"
$ws.Range("BA1").Value = "Outbound Premium Price"

# --- Cells whose literal text would otherwise be auto-detected as a number
#     by Excel (ages like "15", prices like "$601", and the blank
#     Adult-1-Assistance cells) -- force the Text number format first so
#     the value round-trips as a string, matching the rest of the sheet. ---
$forceTextCells = "W19,AS19,Z19,W20,AY20,BA20,Z20,W21,AY21,BA21,Z21"
foreach ($ref in $forceTextCells -split ",") {
    $ws.Range($ref).NumberFormat = "@"
}

# --- Row 19 ---
$ws.Range("A19").Value = "MEKTX"
$ws.Range("C19").Value = "5/3/2025, 13:28:19"
$ws.Range("D19").Value = "SCL"
$ws.Range("E19").Value = "BCN"
$ws.Range("F19").Value = "One Way"
$ws.Range("G19").Value = "Economy"
$ws.Range("H19").Value = "OUTBOUND"
$ws.Range("I19").Value = "Economy"
$ws.Range("J19").Value = "Light"
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = "EN"
$ws.Range("O19").Value = "Juan"
$ws.Range("P19").Value = "Perez"
$ws.Range("Q19").Value = "sofiainkoova@gmail.com"
$ws.Range("R19").Value = "+93 791234567"
$ws.Range("S19").Value = "2222 4000 7000 0005"
$ws.Range("T19").Value = "Juan Pablo Antonio Maximiliano Alejandro Fernández de Córdoba Sánchez"
$ws.Range("U19").Value = "Juan"
$ws.Range("V19").Value = "Perez"
$ws.Range("W19").Value = "15"
$ws.Range("X19").Value = "Afghanistan"
$ws.Range("Y19").Value = "Female"
$ws.Range("Z19").Value = ""
$ws.Range("AC19").Value = "10:25"
$ws.Range("AE19").Value = "06:10"
$ws.Range("AF19").Value = "13h 45m"
$ws.Range("AL19").Value = "Tue, Jul 8"
$ws.Range("AM19").Value = "SCL"
$ws.Range("AN19").Value = "Santiago"
$ws.Range("AO19").Value = "BCN"
$ws.Range("AP19").Value = "Barcelona"
$ws.Range("AQ19").Value = "Yes"
$ws.Range("AR19").Value = "IB"
$ws.Range("AS19").Value = "`$601"

# --- Row 20 ---
$ws.Range("A20").Value = "MEKX4"
$ws.Range("C20").Value = "5/3/2025, 13:38:20"
$ws.Range("D20").Value = "SCL"
$ws.Range("E20").Value = "BCN"
$ws.Range("F20").Value = "One Way"
$ws.Range("G20").Value = "Economy"
$ws.Range("H20").Value = "Light"
$ws.Range("I20").Value = "Economy"
$ws.Range("J20").Value = "Light"
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = "EN"
$ws.Range("O20").Value = "Juan"
$ws.Range("P20").Value = "Perez"
$ws.Range("Q20").Value = "sofiainkoova@gmail.com"
$ws.Range("R20").Value = "+93 791234567"
$ws.Range("S20").Value = "2222 4000 7000 0005"
$ws.Range("T20").Value = "Juan Pablo Antonio Maximiliano Alejandro Fernández de Córdoba Sánchez"
$ws.Range("U20").Value = "Juan"
$ws.Range("V20").Value = "Perez"
$ws.Range("W20").Value = "15"
$ws.Range("X20").Value = "Afghanistan"
$ws.Range("Y20").Value = "Female"
$ws.Range("Z20").Value = ""
$ws.Range("AC20").Value = "10:25"
$ws.Range("AE20").Value = "06:10"
$ws.Range("AF20").Value = "13h 45m"
$ws.Range("AL20").Value = "Tue, Jul 29"
$ws.Range("AM20").Value = "SCL"
$ws.Range("AN20").Value = "Santiago"
$ws.Range("AO20").Value = "BCN"
$ws.Range("AP20").Value = "Barcelona"
$ws.Range("AQ20").Value = "Yes"
$ws.Range("AR20").Value = "IB"
$ws.Range("AT20").Value = "OUTBOUND"
$ws.Range("AU20").Value = "Economy"
$ws.Range("AV20").Value = "Light"
$ws.Range("AW20").Value = "Economy Light"
$ws.Range("AX20").Value = "ECONOMY"
$ws.Range("AY20").Value = "`$565"
$ws.Range("AZ20").Value = "PREMIUM"
$ws.Range("BA20").Value = "`$751"

# --- Row 21 ---
$ws.Range("A21").Value = "MEKZG"
$ws.Range("C21").Value = "5/3/2025, 13:46:13"
$ws.Range("D21").Value = "SCL"
$ws.Range("E21").Value = "BCN"
$ws.Range("F21").Value = "One Way"
$ws.Range("G21").Value = "Economy"
$ws.Range("H21").Value = "Light"
$ws.Range("I21").Value = "Economy"
$ws.Range("J21").Value = "Light"
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = "EN"
$ws.Range("O21").Value = "Juan"
$ws.Range("P21").Value = "Perez"
$ws.Range("Q21").Value = "sofiainkoova@gmail.com"
$ws.Range("R21").Value = "+93 791234567"
$ws.Range("S21").Value = "2222 4000 7000 0005"
$ws.Range("T21").Value = "Juan Pablo Antonio Maximiliano Alejandro Fernández de Córdoba Sánchez"
$ws.Range("U21").Value = "Juan"
$ws.Range("V21").Value = "Perez"
$ws.Range("W21").Value = "15"
$ws.Range("X21").Value = "Afghanistan"
$ws.Range("Y21").Value = "Female"
$ws.Range("Z21").Value = ""
$ws.Range("AC21").Value = "10:25"
$ws.Range("AE21").Value = "06:10"
$ws.Range("AF21").Value = "13h 45m"
$ws.Range("AL21").Value = "Tue, Jul 8"
$ws.Range("AM21").Value = "SCL"
$ws.Range("AN21").Value = "Santiago"
$ws.Range("AO21").Value = "BCN"
$ws.Range("AP21").Value = "Barcelona"
$ws.Range("AQ21").Value = "Yes"
$ws.Range("AR21").Value = "IB"
$ws.Range("AT21").Value = "OUTBOUND"
$ws.Range("AU21").Value = "Economy"
$ws.Range("AV21").Value = "Light"
$ws.Range("AW21").Value = "Economy Light"
$ws.Range("AX21").Value = "ECONOMY"
$ws.Range("AY21").Value = "`$601"
$ws.Range("AZ21").Value = "PREMIUM"
$ws.Range("BA21").Value = "`$751"
